$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-21 Thursday" "2025-08-22 Friday"

Replace-Text "679÷6=" "601÷7="
Replace-Text "217÷2=" "703÷6="
Replace-Text "425÷6=" "931÷2="
Replace-Text "398÷7=" "842÷3="
Replace-Text "591÷2=" "235÷5="
Replace-Text "639÷7=" "163÷9="
Replace-Text "957÷5=" "174÷2="
Replace-Text "572÷8=" "137÷3="
Replace-Text "484÷8=" "994÷4="
Replace-Text "276÷8=" "141÷7="
Replace-Text "215÷2=" "738÷8="
Replace-Text "499÷9=" "956÷4="
Replace-Text "702÷8=" "406÷7="
Replace-Text "269÷7=" "970÷5="
Replace-Text "347÷6=" "973÷9="
Replace-Text "522÷4=" "784÷2="
Replace-Text "433÷9=" "328÷4="
Replace-Text "913÷6=" "417÷3="
Replace-Text "225÷8=" "822÷5="
Replace-Text "355÷9=" "669÷9="
Replace-Text "301÷7=" "763÷2="
Replace-Text "973÷3=" "699÷7="
Replace-Text "769÷4=" "473÷6="
Replace-Text "138÷4=" "641÷7="
Replace-Text "787÷6=" "974÷3="
